$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Sph"
$ws.Range("C2").Value = 0.5362
$ws.Range("D2").Value = 0.9239000000000001
$ws.Range("E2").Value = 0.25
$ws.Range("F2").Value = 0.5803658404589241
$ws.Range("G2").Value = 1743.5906
$ws.Range("H2").Value = 0.2931067733833046
